$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D to hold the two newest quarters
# (2018-12-31 and 2018-09-30). This shifts the existing quarter columns
# (old D:K) right to F:M.
$ws.Range("D:E").Insert()

# Re-apply the correct number formats (date format on the "Period Ending" rows,
# plain number format on all other data rows) to the newly inserted D:E cells by
# copying formats from column F (the original D column, now shifted right).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new quarter values.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value = 142900
$ws.Range("E8").Value = 147700
$ws.Range("D9").Value = 105800
$ws.Range("E9").Value = 103900
$ws.Range("D10").Value = 37100
$ws.Range("E10").Value = 43800
$ws.Range("D12").Value = 1100
$ws.Range("E12").Value = 900
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 3900
$ws.Range("E14").Value = -200
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 138800
$ws.Range("E17").Value = 133700
$ws.Range("D18").Value = 4100
$ws.Range("E18").Value = 14000
$ws.Range("D20").Value = -2100
$ws.Range("E20").Value = -1600
$ws.Range("D21").Value = 9000
$ws.Range("E21").Value = 19600
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 1900
$ws.Range("E23").Value = 12500
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 1900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1500
$ws.Range("E26").Value = 10600
$ws.Range("D27").Value = 1400
$ws.Range("E27").Value = 10500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 2100
$ws.Range("E32").Value = 1600
$ws.Range("D33").Value = 1400
$ws.Range("E33").Value = 10500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 1400
$ws.Range("E35").Value = 10500
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value = 93600
$ws.Range("E41").Value = 85300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 98100
$ws.Range("E43").Value = 123400
$ws.Range("D44").Value = 158500
$ws.Range("E44").Value = 158800
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 350100
$ws.Range("E46").Value = 367500
$ws.Range("D47").Value = 5400
$ws.Range("E47").Value = 5600
$ws.Range("D48").Value = 213300
$ws.Range("E48").Value = 210200
$ws.Range("D49").Value = 35300
$ws.Range("E49").Value = 36200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 12800
$ws.Range("E52").Value = 12900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 616900
$ws.Range("E54").Value = 632400
$ws.Range("D57").Value = 87700
$ws.Range("E57").Value = 82500
$ws.Range("D58").Value = 10500
$ws.Range("E58").Value = 5800
$ws.Range("D59").Value = 13100
$ws.Range("E59").Value = 13800
$ws.Range("D60").Value = 111300
$ws.Range("E60").Value = 102000
$ws.Range("D61").Value = 7100
$ws.Range("E61").Value = 7500
$ws.Range("D62").Value = 32100
$ws.Range("E62").Value = 31200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 144800
$ws.Range("E66").Value = 156000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 360700
$ws.Range("E72").Value = 364400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 472100
$ws.Range("E76").Value = 476400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value = 1400
$ws.Range("E81").Value = 10500
$ws.Range("D83").Value = 7100
$ws.Range("E83").Value = 7200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 36200
$ws.Range("E89").Value = -1200
$ws.Range("D91").Value = -7300
$ws.Range("E91").Value = -5900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -7200
$ws.Range("E94").Value = -6000
$ws.Range("D96").Value = -5200
$ws.Range("E96").Value = -5200
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -20600
$ws.Range("E100").Value = -12000
$ws.Range("D101").Value = -100
$ws.Range("E101").Value = -100
$ws.Range("D102").Value = 8300
$ws.Range("E102").Value = -19400

# Re-fit all column widths now that new data has been added, mirroring what Excel
# does automatically for "best fit" columns when their content changes.
$ws.Columns.AutoFit() | Out-Null
